# Updated symbol list on Tue Jan 24 06:58:48 UTC 2023 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns for the coin rows whose
# values moved since the last scrape. Values are kept as plain text (as in
# the source workbook) by prefixing numeric-looking strings with a leading
# apostrophe, which forces Excel to store them as text instead of numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force text storage so numeric-looking strings (e.g. "319.06",
    # "4.56%") aren't auto-converted to numbers/percentages.
    if ($value -match '^-?[0-9.]+%?$') {
        $ws.Range($range).Value = "'" + $value
    } else {
        $ws.Range($range).Value = $value
    }
}

$updates = @{
    2  = @{ D = "319.06";     E = "4.56%" }
    3  = @{ D = "36.12";      E = "0.23%" }
    4  = @{ D = "5.133";      E = "0.15%" }
    5  = @{ D = "0.08209";    E = "4.34%" }
    6  = @{ D = "2.146";      E = "-1.68%" }
    7  = @{ D = "8.007";      E = "0.89%" }
    8  = @{ D = "4.135" }
    9  = @{ D = "0.9274";     E = "0.90%" }
    10 = @{ D = "0.09977";    E = "2.95%" }
    11 = @{ D = "0.1888";     E = "0.90%" }
    12 = @{ D = "0.09209";    E = "6.31%" }
    13 = @{ D = "0.03617";    E = "3.85%" }
    14 = @{ D = "0.09917";    E = "-0.28%" }
    15 = @{ D = "0.001437";   E = "0.61%" }
    16 = @{ D = "0.005666";   E = "-0.77%" }
    17 = @{ D = "3.468";      E = "0.23%" }
    19 = @{ D = "0.3376";     E = "-1.48%" }
    20 = @{ D = "0.1301";     E = "0.11%" }
    21 = @{ D = "5.058";      E = "4.42%" }
    22 = @{ D = "0.2189";     E = "-0.51%" }
    23 = @{ D = "0.04596" }
    24 = @{ D = "0.001247";   E = "1.31%" }
    25 = @{ D = "0.004734";   E = "-6.99%" }
    26 = @{ D = "0.0001251";  E = "-10.63%" }
    27 = @{ D = "0.0004502";  E = "-5.23%" }
    39 = @{ D = "0.02004";    E = "8.12%" }
    40 = @{ D = "0.04978";    E = "4.35%" }
    41 = @{ D = "0.007766";   E = "0.14%" }
    42 = @{ D = "0.1399";     E = "-0.12%" }
    43 = @{ D = "0.007670";   E = "-0.83%" }
    44 = @{ D = "0.002098";   E = "-5.92%" }
    45 = @{ D = "0.01187";    E = "4.88%" }
    46 = @{ D = "0.00006458"; E = "0.55%" }
    47 = @{ E = "0.09%" }
    48 = @{ D = "29.18";      E = "19.02%" }
    49 = @{ D = "0.001901";   E = "-4.97%" }
    50 = @{ E = "0.09%" }
    51 = @{ E = "0.09%" }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        Set-TextValue "$col$row" $cols[$col]
    }
}
